$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.695.95"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.793.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.791.58"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.718"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000351"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.11"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.382.60"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +20.84%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.825.03"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.137"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.782.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "403.73"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.15"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.99"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.19"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +26.53%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "724.94"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.34"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.120"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.60"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.31"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +24.69%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = [string]::Concat("0.0", [char]0x2083, "0757")
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +15.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0448"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.311"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.20%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.07"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.77"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.65%  "
